# Hortaliza, Vega Modelo de Temuco - Cebollín
# Weekly fruit/vegetable price data refresh: two new daily records are
# inserted near the top of the data block (row 505), pushing all the
# existing records down by two rows (old 505:583 -> new 507:585).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 505, shifting rows 505:583 down to 507:585.
$ws.Rows("505:506").Insert()

# New record in row 505
$ws.Cells.Item(505, 1).Value  = 10
$ws.Cells.Item(505, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(505, 3).Value  = "La Araucanía"
$ws.Cells.Item(505, 4).Value  = 45015
$ws.Cells.Item(505, 5).Value  = 9
$ws.Cells.Item(505, 6).Value  = 100112037
$ws.Cells.Item(505, 7).Value  = "Cebollín"
$ws.Cells.Item(505, 8).Value  = "Sin especificar"
$ws.Cells.Item(505, 9).Value  = "Primera"
$ws.Cells.Item(505, 10).Value = 170
$ws.Cells.Item(505, 11).Value = 7000
$ws.Cells.Item(505, 12).Value = 8000
$ws.Cells.Item(505, 13).Value = 7471
$ws.Cells.Item(505, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(505, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(505, 16).Value = 623
$ws.Cells.Item(505, 17).Value = 12
$ws.Cells.Item(505, 18).Value = "Hortaliza"

# New record in row 506
$ws.Cells.Item(506, 1).Value  = 10
$ws.Cells.Item(506, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(506, 3).Value  = "La Araucanía"
$ws.Cells.Item(506, 4).Value  = 45015
$ws.Cells.Item(506, 5).Value  = 9
$ws.Cells.Item(506, 6).Value  = 100112037
$ws.Cells.Item(506, 7).Value  = "Cebollín"
$ws.Cells.Item(506, 8).Value  = "Sin especificar"
$ws.Cells.Item(506, 9).Value  = "Primera"
$ws.Cells.Item(506, 10).Value = 50
$ws.Cells.Item(506, 11).Value = 8000
$ws.Cells.Item(506, 12).Value = 8000
$ws.Cells.Item(506, 13).Value = 8000
$ws.Cells.Item(506, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(506, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(506, 16).Value = 667
$ws.Cells.Item(506, 17).Value = 12
$ws.Cells.Item(506, 18).Value = "Hortaliza"
